$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of gene data (row 29)
$ws.Range("A29").Value = "Cre02.g095137"
$ws.Range("B29").Value = "PFR1"
$ws.Range("D29").Value = "Pyruvate ferredoxin oxidoreductase"
$ws.Range("G29").Value = "TCA"

# Match the alignment style already used by the rest of column G (left/center)
# by copying the format from the cell above rather than creating a new style.
$ws.Range("G28").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("D32").Select()
